$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing changed values in column B
$ws.Range("B8").Value = -0.957382811
$ws.Range("B9").Value = -3.098252891
$ws.Range("B13").Value = -0.09751874099999999
$ws.Range("B16").Value = 0.289898428
$ws.Range("B18").Value = 0.08153150799999997
$ws.Range("B19").Value = 0.754887929
$ws.Range("B20").Value = 0.427370554
$ws.Range("B21").Value = 0.615441003
$ws.Range("B22").Value = 0.238543425
$ws.Range("B23").Value = -0.255775563

# Add new rows 24 and 25
$ws.Range("A24").Value = "2025-07-01_diff"
$ws.Range("B24").Value = 0.407980578
$ws.Range("A25").Value = "2025-10-01_diff"

# Apply same style (style index 1 in original) to the new A cells as other date-label cells
$ws.Range("A23").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
